$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2431
$ws.Range("F7").Value = 77
$ws.Range("F9").Value = 1697
$ws.Range("F10").Value = 1697
$ws.Range("F11").Value = 1401
$ws.Range("F12").Value = 79
$ws.Range("F16").Value = 903
$ws.Range("F17").Value = 225
$ws.Range("F18").Value = 191
$ws.Range("F19").Value = 243
$ws.Range("F20").Value = 7533
$ws.Range("F21").Value = 8530
$ws.Range("F24").Value = 421
$ws.Range("F30").Value = 14
$ws.Range("F33").Value = 1507
$ws.Range("F34").Value = 12
$ws.Range("F35").Value = 262
$ws.Range("F42").Value = 1374
$ws.Range("F43").Value = 371
$ws.Range("F45").Value = 222
$ws.Range("F48").Value = 192
$ws.Range("F49").Value = 38
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 22
$ws.Range("F4").Value = 10
$ws.Range("F15").Value = 20
$ws.Range("F17").Value = 1
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 191
$ws.Range("F4").Value = 303
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 191
$ws.Range("F5").Value = 22
$ws.Range("F6").Value = 303
$ws.Range("F10").Value = 2431
$ws.Range("F12").Value = 77
$ws.Range("F14").Value = 1697
$ws.Range("F15").Value = 1697
$ws.Range("F16").Value = 1401
$ws.Range("F17").Value = 79
$ws.Range("F20").Value = 225
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 191
$ws.Range("F24").Value = 243
$ws.Range("F25").Value = 7533
$ws.Range("F26").Value = 8530
$ws.Range("F28").Value = 421
$ws.Range("F32").Value = 14
$ws.Range("F34").Value = 1507
$ws.Range("F35").Value = 12
$ws.Range("F36").Value = 262
$ws.Range("F43").Value = 1374
$ws.Range("F44").Value = 371
$ws.Range("F46").Value = 222
$ws.Range("F48").Value = 192
$ws.Range("F50").Value = 38
